$p = $ppt.ActivePresentation
$m = $p.SlideMaster
$t = $m.Theme
$tv = $t.ThemeVariants
Write-Output $tv
$m2 = Get-Member -InputObject $tv
Write-Output ($m2 | Out-String)
